$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(39.06, 0.04000000000000001, 3.029262542724609),
    @(39.06, 0.04000000000000001, 3.055202484130859),
    @(39.06, 0.04000000000000001, 3.826434135437012),
    @(37.37, 0.36, 22.244460105896),
    @(39.06, 0.04000000000000001, 3.323311328887939),
    @(37.37, 0.36, 20.74647331237793),
    @(39.06, 0.04000000000000001, 2.801284313201904),
    @(37.37, 0.36, 20.7770311832428),
    @(39.06, 0.04000000000000001, 3.219979524612427),
    @(37.37, 0.36, 20.27344155311584),
    @(39.06, 0.04000000000000001, 2.994482517242432),
    @(39.06, 0.04000000000000001, 2.753181219100952),
    @(39.06, 0.04000000000000001, 2.686118364334106),
    @(37.37, 0.36, 20.11983466148376),
    @(39.06, 0.04000000000000001, 3.07470440864563),
    @(39.06, 0.04000000000000001, 3.950433254241943),
    @(39.06, 0.04000000000000001, 3.116287708282471)
)

$startRow = 4
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
